$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 624.3077
$ws.Range("I38").Value = 226
$ws.Range("J38").Value = 1261.6
$ws.Range("K38").Value = 678
$ws.Range("L38").Value = 3784.8
$ws.Range("M38").Value = -306
$ws.Range("N38").Value = -4528.799999999999
$ws.Range("H58").Value = 2472.1785
$ws.Range("I58").Value = 471.66666
$ws.Range("J58").Value = 2712.24
$ws.Range("K58").Value = 1414.99998
$ws.Range("L58").Value = 8136.719999999999
$ws.Range("M58").Value = -1264.99998
$ws.Range("N58").Value = -8436.719999999999
$ws.Range("H87").Value = 39539.6
$ws.Range("J87").Value = 39539.6
$ws.Range("L87").Value = 39539.6
$ws.Range("N87").Value = -42035.6
$ws.Range("H90").Value = 39539.6
$ws.Range("J90").Value = 39539.6
$ws.Range("L90").Value = 118618.8
$ws.Range("N90").Value = -131098.8
$ws.Range("H92").Value = 962010.4399999999
$ws.Range("I92").Value = 1282347.2
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 1282347.2
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = -1281099.2
$ws.Range("N92").Value = -3496
$ws.Range("H110").Value = 20711.555
$ws.Range("J110").Value = 20711.555
$ws.Range("L110").Value = 20711.555
$ws.Range("N110").Value = -28891.555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5047.055
$ws.Range("I32").Value = 4678.024
$ws.Range("J32").Value = 9475.429
$ws.Range("K32").Value = 4678.024
$ws.Range("L32").Value = 9475.429
$ws.Range("M32").Value = -4391.024
$ws.Range("N32").Value = -10049.429
$ws.Range("H112").Value = 23462.334
$ws.Range("J112").Value = 23462.334
$ws.Range("L112").Value = 23462.334
$ws.Range("N112").Value = -26416.334
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H132").Value = 3848.6619
$ws.Range("I132").Value = 2418.5435
$ws.Range("J132").Value = 6480.08
$ws.Range("K132").Value = 7255.630500000001
$ws.Range("L132").Value = 19440.24
$ws.Range("M132").Value = -4725.630500000001
$ws.Range("N132").Value = -24500.24

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5167.702
$ws.Range("I134").Value = 2236.4167
$ws.Range("J134").Value = 8226.434999999999
$ws.Range("K134").Value = 6709.250100000001
$ws.Range("L134").Value = 24679.305
$ws.Range("M134").Value = -4174.250100000001
$ws.Range("N134").Value = -29749.305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7938752
$ws.Range("I31").Value = 1603.6666
$ws.Range("J31").Value = 18521616
$ws.Range("K31").Value = 1603.6666
$ws.Range("L31").Value = 18521616
$ws.Range("M31").Value = -1308.6666
$ws.Range("N31").Value = -18522206
$ws.Range("H34").Value = 7938752
$ws.Range("I34").Value = 1603.6666
$ws.Range("J34").Value = 18521616
$ws.Range("K34").Value = 1603.6666
$ws.Range("L34").Value = 18521616
$ws.Range("M34").Value = -1401.6666
$ws.Range("N34").Value = -18522020
$ws.Range("H118").Value = 61370
$ws.Range("J118").Value = 61370
$ws.Range("L118").Value = 61370
$ws.Range("N118").Value = -64684
$ws.Range("H132").Value = 2427.25
$ws.Range("I132").Value = 1814.3572
$ws.Range("J132").Value = 3285.3
$ws.Range("K132").Value = 5443.071599999999
$ws.Range("L132").Value = 9855.900000000001
$ws.Range("M132").Value = -2913.071599999999
$ws.Range("N132").Value = -14915.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2386.3635
$ws.Range("I98").Value = 603.6667
$ws.Range("J98").Value = 4525.6
$ws.Range("K98").Value = 1811.0001
$ws.Range("L98").Value = 13576.8
$ws.Range("M98").Value = -313.0001
$ws.Range("N98").Value = -16572.8
$ws.Range("H132").Value = 2631
$ws.Range("I132").Value = 787.5
$ws.Range("K132").Value = 7087.5
$ws.Range("M132").Value = -4557.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 6000
$ws.Range("J47").Value = 6000
$ws.Range("L47").Value = 6000
$ws.Range("N47").Value = -7136
$ws.Range("H55").Value = 5990
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5990
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 5990
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -6644
$ws.Range("H132").Value = 1304811.6
$ws.Range("I132").Value = 3790097.8
$ws.Range("J132").Value = 2995.1428
$ws.Range("K132").Value = 11370293.4
$ws.Range("L132").Value = 8985.428400000001
$ws.Range("M132").Value = -11367763.4
$ws.Range("N132").Value = -14045.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2135.7144
$ws.Range("I7").Value = 2241.6667
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 2241.6667
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -2129.6667
$ws.Range("N7").Value = -1724
$ws.Range("H101").Value = 15172.2
$ws.Range("J101").Value = 15172.2
$ws.Range("L101").Value = 15172.2
$ws.Range("N101").Value = -21662.2
$ws.Range("H110").Value = 14988.8
$ws.Range("J110").Value = 14988.8
$ws.Range("L110").Value = 14988.8
$ws.Range("N110").Value = -23168.8
$ws.Range("H126").Value = 2135.7144
$ws.Range("I126").Value = 2241.6667
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 6725.000100000001
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -4255.000100000001
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 38464610
$ws.Range("I132").Value = 45457436
$ws.Range("J132").Value = 4057
$ws.Range("K132").Value = 136372308
$ws.Range("L132").Value = 12171
$ws.Range("M132").Value = -136369778
$ws.Range("N132").Value = -17231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9999.75
$ws.Range("I54").Value = 3000
$ws.Range("J54").Value = 12333
$ws.Range("K54").Value = 3000
$ws.Range("L54").Value = 12333
$ws.Range("M54").Value = -2480
$ws.Range("N54").Value = -13373
$ws.Range("H81").Value = 2257.889
$ws.Range("I81").Value = 2667
$ws.Range("J81").Value = 1439.6666
$ws.Range("K81").Value = 5334
$ws.Range("L81").Value = 2879.3332
$ws.Range("M81").Value = -4273
$ws.Range("N81").Value = -5001.3332
$ws.Range("H84").Value = 2257.889
$ws.Range("I84").Value = 2667
$ws.Range("J84").Value = 1439.6666
$ws.Range("K84").Value = 26670
$ws.Range("L84").Value = 14396.666
$ws.Range("M84").Value = -21366
$ws.Range("N84").Value = -25004.666
$ws.Range("H103").Value = 21447.5
$ws.Range("J103").Value = 21447.5
$ws.Range("L103").Value = 21447.5
$ws.Range("N103").Value = -23791.5
